$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "H2-M3"
$ws.Range("C2").Value = "Klrd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.366593333333333
$ws.Range("H2").Value = 13.09978
$ws.Range("I2").Value = 0.1610144701463773
$ws.Range("J2").Value = 0.1610144701463773
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04040533333333333
$ws.Range("N2").Value = 0.121216
$ws.Range("O2").Value = 0.009158328487187175
$ws.Range("P2").Value = 0.009158328487187173
$ws.Range("Q2").Value = 0.1764336591644445
$ws.Range("R2").Value = 1.58790293248
$ws.Range("S2").Value = 0.001474623408790916
$ws.Range("T2").Value = 0.001474623408790915

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "H2-M3"
$ws.Range("C3").Value = "Klrd1"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.366593333333333
$ws.Range("H3").Value = 13.09978
$ws.Range("I3").Value = 0.1610144701463773
$ws.Range("J3").Value = 0.1610144701463773
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.697694
$ws.Range("N3").Value = 8.093081999999999
$ws.Range("O3").Value = 0.6114630364782021
$ws.Range("P3").Value = 0.611463036478202
$ws.Range("Q3").Value = 11.77973263577333
$ws.Range("R3").Value = 106.01759372196
$ws.Range("S3").Value = 0.09845439683263267
$ws.Range("T3").Value = 0.09845439683263264

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "H2-M3"
$ws.Range("C4").Value = "Klrd1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.366593333333333
$ws.Range("H4").Value = 13.09978
$ws.Range("I4").Value = 0.1610144701463773
$ws.Range("J4").Value = 0.1610144701463773
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.673768333333333
$ws.Range("N4").Value = 5.021305
$ws.Range("O4").Value = 0.3793786350346109
$ws.Range("P4").Value = 0.3793786350346108
$ws.Range("Q4").Value = 7.308665645877777
$ws.Range("R4").Value = 65.7779908129
$ws.Range("S4").Value = 0.06108544990495371
$ws.Range("T4").Value = 0.0610854499049537

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "H2-M3"
$ws.Range("C5").Value = "Klrd1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.732182666666667
$ws.Range("H5").Value = 11.196548
$ws.Range("I5").Value = 0.1376211084223155
$ws.Range("J5").Value = 0.1376211084223155
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04040533333333333
$ws.Range("N5").Value = 0.121216
$ws.Range("O5").Value = 0.009158328487187175
$ws.Range("P5").Value = 0.009158328487187173
$ws.Range("Q5").Value = 0.1508000847075556
$ws.Range("R5").Value = 1.357200762368
$ws.Range("S5").Value = 0.001260379317702367
$ws.Range("T5").Value = 0.001260379317702367

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "H2-M3"
$ws.Range("C6").Value = "Klrd1"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.732182666666667
$ws.Range("H6").Value = 11.196548
$ws.Range("I6").Value = 0.1376211084223155
$ws.Range("J6").Value = 0.1376211084223155
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.697694
$ws.Range("N6").Value = 8.093081999999999
$ws.Range("O6").Value = 0.6114630364782021
$ws.Range("P6").Value = 0.611463036478202
$ws.Range("Q6").Value = 10.06828678677067
$ws.Range("R6").Value = 90.61458108093599
$ws.Range("S6").Value = 0.08415022083940492
$ws.Range("T6").Value = 0.0841502208394049

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "H2-M3"
$ws.Range("C7").Value = "Klrd1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.732182666666667
$ws.Range("H7").Value = 11.196548
$ws.Range("I7").Value = 0.1376211084223155
$ws.Range("J7").Value = 0.1376211084223155
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.673768333333333
$ws.Range("N7").Value = 5.021305
$ws.Range("O7").Value = 0.3793786350346109
$ws.Range("P7").Value = 0.3793786350346108
$ws.Range("Q7").Value = 6.246809161682222
$ws.Range("R7").Value = 56.22128245514
$ws.Range("S7").Value = 0.05221050826520825
$ws.Range("T7").Value = 0.05221050826520824

$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "H2-M3"
$ws.Range("C8").Value = "Klrd1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.955599666666666
$ws.Range("H8").Value = 26.866799
$ws.Range("I8").Value = 0.3302302332950797
$ws.Range("J8").Value = 0.3302302332950797
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04040533333333333
$ws.Range("N8").Value = 0.121216
$ws.Range("O8").Value = 0.009158328487187175
$ws.Range("P8").Value = 0.009158328487187173
$ws.Range("Q8").Value = 0.3618539897315555
$ws.Range("R8").Value = 3.256685907584
$ws.Range("S8").Value = 0.003024356952916795
$ws.Range("T8").Value = 0.003024356952916794

$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "H2-M3"
$ws.Range("C9").Value = "Klrd1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.955599666666666
$ws.Range("H9").Value = 26.866799
$ws.Range("I9").Value = 0.3302302332950797
$ws.Range("J9").Value = 0.3302302332950797
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.697694
$ws.Range("N9").Value = 8.093081999999999
$ws.Range("O9").Value = 0.6114630364782021
$ws.Range("P9").Value = 0.611463036478202
$ws.Range("Q9").Value = 24.15946748716866
$ws.Range("R9").Value = 217.435207384518
$ws.Range("S9").Value = 0.2019235811875145
$ws.Range("T9").Value = 0.2019235811875144

$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "H2-M3"
$ws.Range("C10").Value = "Klrd1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.955599666666666
$ws.Range("H10").Value = 26.866799
$ws.Range("I10").Value = 0.3302302332950797
$ws.Range("J10").Value = 0.3302302332950797
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.673768333333333
$ws.Range("N10").Value = 5.021305
$ws.Range("O10").Value = 0.3793786350346109
$ws.Range("P10").Value = 0.3793786350346108
$ws.Range("Q10").Value = 14.98959912807722
$ws.Range("R10").Value = 134.906392152695
$ws.Range("S10").Value = 0.1252822951546484
$ws.Range("T10").Value = 0.1252822951546484

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "H2-M3"
$ws.Range("C11").Value = "Klrd1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.764211666666668
$ws.Range("H11").Value = 29.292635
$ws.Range("I11").Value = 0.3600471232124682
$ws.Range("J11").Value = 0.3600471232124682
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04040533333333333
$ws.Range("N11").Value = 0.121216
$ws.Range("O11").Value = 0.009158328487187175
$ws.Range("P11").Value = 0.009158328487187173
$ws.Range("Q11").Value = 0.394526227128889
$ws.Range("R11").Value = 3.550736044160001
$ws.Range("S11").Value = 0.003297429825246538
$ws.Range("T11").Value = 0.003297429825246537

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "H2-M3"
$ws.Range("C12").Value = "Klrd1"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.764211666666668
$ws.Range("H12").Value = 29.292635
$ws.Range("I12").Value = 0.3600471232124682
$ws.Range("J12").Value = 0.3600471232124682
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.697694
$ws.Range("N12").Value = 8.093081999999999
$ws.Range("O12").Value = 0.6114630364782021
$ws.Range("P12").Value = 0.611463036478202
$ws.Range("Q12").Value = 26.34085522789667
$ws.Range("R12").Value = 237.06769705107
$ws.Range("S12").Value = 0.2201555072347371
$ws.Range("T12").Value = 0.2201555072347371

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "H2-M3"
$ws.Range("C13").Value = "Klrd1"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.764211666666668
$ws.Range("H13").Value = 29.292635
$ws.Range("I13").Value = 0.3600471232124682
$ws.Range("J13").Value = 0.3600471232124682
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.673768333333333
$ws.Range("N13").Value = 5.021305
$ws.Range("O13").Value = 0.3793786350346109
$ws.Range("P13").Value = 0.3793786350346108
$ws.Range("Q13").Value = 16.34302828763056
$ws.Range("R13").Value = 147.087254588675
$ws.Range("S13").Value = 0.1365941861524845
$ws.Range("T13").Value = 0.1365941861524845

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "H2-M3"
$ws.Range("C14").Value = "Klrd1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.300673
$ws.Range("H14").Value = 0.9020190000000001
$ws.Range("I14").Value = 0.01108706492375941
$ws.Range("J14").Value = 0.01108706492375941
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.04040533333333333
$ws.Range("N14").Value = 0.121216
$ws.Range("O14").Value = 0.009158328487187175
$ws.Range("P14").Value = 0.009158328487187173
$ws.Range("Q14").Value = 0.01214879278933334
$ws.Range("R14").Value = 0.109339135104
$ws.Range("S14").Value = 0.0001015389825305595
$ws.Range("T14").Value = 0.0001015389825305595

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "H2-M3"
$ws.Range("C15").Value = "Klrd1"
$ws.Range("D15").Value = "M1"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.300673
$ws.Range("H15").Value = 0.9020190000000001
$ws.Range("I15").Value = 0.01108706492375941
$ws.Range("J15").Value = 0.01108706492375941
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.697694
$ws.Range("N15").Value = 8.093081999999999
$ws.Range("O15").Value = 0.6114630364782021
$ws.Range("P15").Value = 0.611463036478202
$ws.Range("Q15").Value = 0.811123748062
$ws.Range("R15").Value = 7.300113732558001
$ws.Range("S15").Value = 0.006779330383912897
$ws.Range("T15").Value = 0.006779330383912896

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "H2-M3"
$ws.Range("C16").Value = "Klrd1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.300673
$ws.Range("H16").Value = 0.9020190000000001
$ws.Range("I16").Value = 0.01108706492375941
$ws.Range("J16").Value = 0.01108706492375941
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.673768333333333
$ws.Range("N16").Value = 5.021305
$ws.Range("O16").Value = 0.3793786350346109
$ws.Range("P16").Value = 0.3793786350346108
$ws.Range("Q16").Value = 0.5032569460883334
$ws.Range("R16").Value = 4.529312514795
$ws.Range("S16").Value = 0.004206195557315958
$ws.Range("T16").Value = 0.004206195557315958
